$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Widen column A slightly to fit the new, longer item names ---
$ws.Columns.Item(1).ColumnWidth = 40.83

# --- Row 5: Switching Power Supply 24V 20A ---
$ws.Range("A5").Value2 = "Switching Power Supply 24V 20A"
$ws.Range("B5").Value2 = 1
$ws.Range("C5").Value2 = 37.99

$ws.Range("E3").Copy()
$ws.Range("E5").PasteSpecial(-4122)
$ws.Range("E5").Value2 = 43532

$ws.Range("F5").Value2 = "Amazon"

$ws.Range("G5").Value2 = "Click Here"
$ws.Hyperlinks.Add($ws.Range("G5"), "https://www.amazon.co.uk/dp/B01MZZTLDL", [Type]::Missing, [Type]::Missing, "Click Here")
$ws.Range("G5").Style = "Hyperlink"

# --- Row 6: 3PCS TB6600 4A 9-42V Stepper Motor Driver ---
$ws.Range("A6").Value2 = "3PCS TB6600 4A 9-42V Stepper Motor Driver"
$ws.Range("B6").Value2 = 1
$ws.Range("C6").Value2 = 22.99

$ws.Range("E3").Copy()
$ws.Range("E6").PasteSpecial(-4122)
$ws.Range("E6").Value2 = 43532

$ws.Range("F6").Value2 = "Amazon"

$ws.Range("G6").Value2 = "Click Here"
$ws.Hyperlinks.Add($ws.Range("G6"), "https://www.amazon.co.uk/dp/B01J3UKG1C", [Type]::Missing, [Type]::Missing, "Click Here")
$ws.Range("G6").Style = "Hyperlink"

# --- Move the active selection, as recorded by Excel on save ---
$ws.Range("F9").Select()
